# Auto-generated Excel COM-interop script
# Applies the "Add data for 2024-09-05" update to violent-crime-full-year.xlsx
# 164 cell updates across 49 worksheets (Citywide Totals, By Neighborhood, and 47 neighborhood sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5467
$ws.Range("K3").Value = 5610
$ws.Range("J4").Value = 1831
$ws.Range("K4").Value = 1168
$ws.Range("K5").Value = 401
$ws.Range("K6").Value = 6237
$ws.Range("J7").Value = 29297
$ws.Range("K7").Value = 18883

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 345
$ws.Range("K6").Value = 425
$ws.Range("K7").Value = 1257

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 221
$ws.Range("K3").Value = 297
$ws.Range("K4").Value = 37
$ws.Range("K6").Value = 237
$ws.Range("K7").Value = 809

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 107
$ws.Range("K3").Value = 113
$ws.Range("K7").Value = 321

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 185
$ws.Range("K6").Value = 186
$ws.Range("K7").Value = 637

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 108
$ws.Range("K7").Value = 431

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 78
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 317

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 231
$ws.Range("K7").Value = 557
$ws.Range("K8").Value = 1257
$ws.Range("K9").Value = 78
$ws.Range("K11").Value = 359
$ws.Range("K15").Value = 190
$ws.Range("K17").Value = 36
$ws.Range("K19").Value = 552
$ws.Range("K20").Value = 436
$ws.Range("K21").Value = 60
$ws.Range("K22").Value = 50
$ws.Range("K24").Value = 57
$ws.Range("K25").Value = 89
$ws.Range("K29").Value = 1015
$ws.Range("K33").Value = 809
$ws.Range("K37").Value = 637
$ws.Range("K40").Value = 44
$ws.Range("K41").Value = 130
$ws.Range("K43").Value = 167
$ws.Range("K44").Value = 164
$ws.Range("K48").Value = 239
$ws.Range("K49").Value = 104
$ws.Range("K51").Value = 236
$ws.Range("K52").Value = 493
$ws.Range("K54").Value = 366
$ws.Range("K56").Value = 20
$ws.Range("K57").Value = 73
$ws.Range("K58").Value = 12
$ws.Range("K63").Value = 52
$ws.Range("K64").Value = 123
$ws.Range("K65").Value = 431
$ws.Range("K66").Value = 61
$ws.Range("K67").Value = 719
$ws.Range("K73").Value = 164
$ws.Range("K76").Value = 262
$ws.Range("K78").Value = 221
$ws.Range("K79").Value = 480
$ws.Range("K84").Value = 144
$ws.Range("K85").Value = 890
$ws.Range("K89").Value = 275
$ws.Range("K91").Value = 211
$ws.Range("K94").Value = 256
$ws.Range("K95").Value = 321
$ws.Range("K97").Value = 149
$ws.Range("K99").Value = 317
$ws.Range("J101").Value = 29297
$ws.Range("K101").Value = 18883

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 207
$ws.Range("K6").Value = 199
$ws.Range("K7").Value = 719

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 366

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 289
$ws.Range("K6").Value = 287
$ws.Range("K7").Value = 1015

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 165
$ws.Range("K6").Value = 174
$ws.Range("K7").Value = 552

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 57
$ws.Range("K6").Value = 137
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 66
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 221

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 54

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 98
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 160
$ws.Range("K6").Value = 121
$ws.Range("K7").Value = 480

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 145
$ws.Range("K7").Value = 436

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 188
$ws.Range("K7").Value = 557

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 50
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 38
$ws.Range("K3").Value = 35

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 67
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 122
$ws.Range("K3").Value = 95
$ws.Range("K7").Value = 359

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 54
$ws.Range("K7").Value = 164
$ws.Range("J4").Value = 17

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 85
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 303
$ws.Range("K6").Value = 218
$ws.Range("K7").Value = 890

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("K6").Value = 7

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 137
$ws.Range("K6").Value = 180
$ws.Range("K7").Value = 493

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 12

